$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "Cotação" header to "Compra R$"
$ws.Range("B1").Value = "Compra R$"

# Add the new "Venda R$" column, matching the header formatting of B1/A1
$ws.Range("C1").Value = "Venda R$"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Fill in the buy/sell (venda) values for each currency row
$ws.Range("C2").Value = 5.8341
$ws.Range("C3").Value = 6.0173
